$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume data scraped on Tue Apr 25 13:25:49 UTC 2023
# Price/Volume columns (D, E) hold numeric-looking text (e.g. "1.000", "27.445.70")
# that Excel would otherwise auto-coerce to a Number/Date on assignment, so the
# cell is pre-formatted as Text ("@") to preserve the literal string content.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.445.70'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.829.79'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -1.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '331.21'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4587'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3823'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07897'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9668'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.87%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.10'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.831.43'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.880'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.094'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '89.66'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06596'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.18'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '27.433.99'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.326'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -3.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.83'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.046.95'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.40'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.073'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.294'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -3.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '118.32'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.93%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09311'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9389'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -4.54%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.81%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.330'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05936'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02178'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.135'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.142'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -4.95%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5784'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -3.38%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1828'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.994'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.265'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5441'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.07%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.875'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06581'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.74%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '110.03'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -33.70%  '
